# Aggiornamento dati fino al 28/06 incluso (righe 270-301)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 270
$endRow = 301

# Replicate the formatting of the last existing data row (269, columns A:D)
# onto the newly appended rows, so column A keeps its date style (s="2")
# and columns B:D stay unstyled, just like the rest of the table.
$ws.Range("A269:D269").Copy()
$ws.Range("A" + $startRow + ":D" + $endRow).PasteSpecial(-4122)

$data = @(
    @(270, 44344, 2, 19, 57.60543310190099),
    @(271, 44345, 0, 15, 45.47797350150077),
    @(272, 44346, 0, 15, 45.47797350150077),
    @(273, 44347, 0, 11, 33.35051390110057),
    @(274, 44348, 0, 7, 21.22305430070036),
    @(275, 44349, 0, 6, 18.19118940060031),
    @(276, 44350, 0, 2, 6.063729800200103),
    @(277, 44351, 1, 1, 3.031864900100051),
    @(278, 44352, 0, 1, 3.031864900100051),
    @(279, 44353, 2, 3, 9.095594700300154),
    @(280, 44354, 0, 3, 9.095594700300154),
    @(281, 44355, 0, 3, 9.095594700300154),
    @(282, 44356, 0, 3, 9.095594700300154),
    @(283, 44357, 0, 3, 9.095594700300154),
    @(284, 44358, 0, 2, 6.063729800200103),
    @(285, 44359, 0, 2, 6.063729800200103),
    @(286, 44360, 0, 0, 0),
    @(287, 44361, 0, 0, 0),
    @(288, 44362, 1, 1, 3.031864900100051),
    @(289, 44363, 0, 1, 3.031864900100051),
    @(290, 44364, 3, 4, 12.12745960040021),
    @(291, 44365, 0, 4, 12.12745960040021),
    @(292, 44366, 1, 5, 15.15932450050026),
    @(293, 44367, 0, 5, 15.15932450050026),
    @(294, 44368, 1, 6, 18.19118940060031),
    @(295, 44369, 0, 5, 15.15932450050026),
    @(296, 44370, 0, 5, 15.15932450050026),
    @(297, 44371, 1, 3, 9.095594700300154),
    @(298, 44372, 0, 3, 9.095594700300154),
    @(299, 44373, 0, 2, 6.063729800200103),
    @(300, 44374, 0, 2, 6.063729800200103),
    @(301, 44375, 1, 2, 6.063729800200103)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Host "Updated rows $startRow to $endRow"
